# Update TPM-derived metrics in Spp1-Itga4 LR-pairs sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = [double]"10.70913233333333"
$ws.Range("H2").Value = [double]"32.127397"
$ws.Range("I2").Value = [double]"0.007451729107954897"
$ws.Range("J2").Value = [double]"0.007451729107954897"
$ws.Range("M2").Value = [double]"0.3331066666666667"
$ws.Range("N2").Value = [double]"0.99932"
$ws.Range("O2").Value = [double]"0.002125805913843485"
$ws.Range("P2").Value = [double]"0.002125805913843485"
$ws.Range("Q2").Value = [double]"3.567283374448889"
$ws.Range("R2").Value = [double]"32.10555037004"
$ws.Range("S2").Value = [double]"1.584092980605015E-05"
$ws.Range("T2").Value = [double]"1.584092980605015E-05"
$ws.Range("G3").Value = [double]"10.70913233333333"
$ws.Range("H3").Value = [double]"32.127397"
$ws.Range("I3").Value = [double]"0.007451729107954897"
$ws.Range("J3").Value = [double]"0.007451729107954897"
$ws.Range("O3").Value = [double]"0.0008775937418887864"
$ws.Range("P3").Value = [double]"0.0008775937418887864"
$ws.Range("Q3").Value = [double]"1.472677041950667"
$ws.Range("R3").Value = [double]"13.254093377556"
$ws.Range("S3").Value = [double]"6.539590831391726E-06"
$ws.Range("T3").Value = [double]"6.539590831391726E-06"
$ws.Range("G4").Value = [double]"10.70913233333333"
$ws.Range("H4").Value = [double]"32.127397"
$ws.Range("I4").Value = [double]"0.007451729107954897"
$ws.Range("J4").Value = [double]"0.007451729107954897"
$ws.Range("M4").Value = [double]"91.40156066666667"
$ws.Range("N4").Value = [double]"274.204682"
$ws.Range("O4").Value = [double]"0.5833025803538128"
$ws.Range("P4").Value = [double]"0.5833025803538128"
$ws.Range("Q4").Value = [double]"978.8314086525284"
$ws.Range("R4").Value = [double]"8809.482677872755"
$ws.Range("S4").Value = [double]"0.004346612816767707"
$ws.Range("T4").Value = [double]"0.004346612816767707"
$ws.Range("G5").Value = [double]"10.70913233333333"
$ws.Range("H5").Value = [double]"32.127397"
$ws.Range("I5").Value = [double]"0.007451729107954897"
$ws.Range("J5").Value = [double]"0.007451729107954897"
$ws.Range("M5").Value = [double]"0.5759770000000001"
$ws.Range("N5").Value = [double]"1.727931"
$ws.Range("O5").Value = [double]"0.00367574544541637"
$ws.Range("P5").Value = [double]"0.00367574544541637"
$ws.Range("Q5").Value = [double]"6.168213913956335"
$ws.Range("R5").Value = [double]"55.513925225607"
$ws.Range("S5").Value = [double]"2.73906593290418E-05"
$ws.Range("T5").Value = [double]"2.73906593290418E-05"
$ws.Range("G6").Value = [double]"10.70913233333333"
$ws.Range("H6").Value = [double]"32.127397"
$ws.Range("I6").Value = [double]"0.007451729107954897"
$ws.Range("J6").Value = [double]"0.007451729107954897"
$ws.Range("M6").Value = [double]"64.24849033333334"
$ws.Range("N6").Value = [double]"192.745471"
$ws.Range("O6").Value = [double]"0.4100182745450386"
$ws.Range("P6").Value = [double]"0.4100182745450385"
$ws.Range("Q6").Value = [double]"688.0455851965543"
$ws.Range("R6").Value = [double]"6192.410266768988"
$ws.Range("S6").Value = [double]"0.003055345111220706"
$ws.Range("T6").Value = [double]"0.003055345111220706"
$ws.Range("I7").Value = [double]"0.03290895798513831"
$ws.Range("J7").Value = [double]"0.03290895798513832"
$ws.Range("M7").Value = [double]"0.3331066666666667"
$ws.Range("N7").Value = [double]"0.99932"
$ws.Range("O7").Value = [double]"0.002125805913843485"
$ws.Range("P7").Value = [double]"0.002125805913843485"
$ws.Range("Q7").Value = [double]"15.75413934002222"
$ws.Range("R7").Value = [double]"141.7872540602"
$ws.Range("S7").Value = [double]"6.995805750323381E-05"
$ws.Range("T7").Value = [double]"6.995805750323382E-05"
$ws.Range("I8").Value = [double]"0.03290895798513831"
$ws.Range("J8").Value = [double]"0.03290895798513832"
$ws.Range("O8").Value = [double]"0.0008775937418887864"
$ws.Range("P8").Value = [double]"0.0008775937418887864"
$ws.Range("R8").Value = [double]"58.53385110678001"
$ws.Range("S8").Value = [double]"2.888069557983839E-05"
$ws.Range("T8").Value = [double]"2.88806955798384E-05"
$ws.Range("I9").Value = [double]"0.03290895798513831"
$ws.Range("J9").Value = [double]"0.03290895798513832"
$ws.Range("M9").Value = [double]"91.40156066666667"
$ws.Range("N9").Value = [double]"274.204682"
$ws.Range("O9").Value = [double]"0.5833025803538128"
$ws.Range("P9").Value = [double]"0.5833025803538128"
$ws.Range("Q9").Value = [double]"4322.798270738585"
$ws.Range("R9").Value = [double]"38905.18443664727"
$ws.Range("S9").Value = [double]"0.01919588010948639"
$ws.Range("T9").Value = [double]"0.01919588010948639"
$ws.Range("I10").Value = [double]"0.03290895798513831"
$ws.Range("J10").Value = [double]"0.03290895798513832"
$ws.Range("M10").Value = [double]"0.5759770000000001"
$ws.Range("N10").Value = [double]"1.727931"
$ws.Range("O10").Value = [double]"0.00367574544541637"
$ws.Range("P10").Value = [double]"0.00367574544541637"
$ws.Range("Q10").Value = [double]"27.24058934469834"
$ws.Range("R10").Value = [double]"245.165304102285"
$ws.Range("S10").Value = [double]"0.0001209649524272708"
$ws.Range("T10").Value = [double]"0.0001209649524272709"
$ws.Range("I11").Value = [double]"0.03290895798513831"
$ws.Range("J11").Value = [double]"0.03290895798513832"
$ws.Range("M11").Value = [double]"64.24849033333334"
$ws.Range("N11").Value = [double]"192.745471"
$ws.Range("O11").Value = [double]"0.4100182745450386"
$ws.Range("P11").Value = [double]"0.4100182745450385"
$ws.Range("Q11").Value = [double]"3038.605258868243"
$ws.Range("R11").Value = [double]"27347.44732981419"
$ws.Range("S11").Value = [double]"0.01349327417014158"
$ws.Range("T11").Value = [double]"0.01349327417014158"
$ws.Range("G12").Value = [double]"411.37678"
$ws.Range("H12").Value = [double]"1234.13034"
$ws.Range("I12").Value = [double]"0.2862480573072345"
$ws.Range("J12").Value = [double]"0.2862480573072345"
$ws.Range("M12").Value = [double]"0.3331066666666667"
$ws.Range("N12").Value = [double]"0.99932"
$ws.Range("O12").Value = [double]"0.002125805913843485"
$ws.Range("P12").Value = [double]"0.002125805913843485"
$ws.Range("Q12").Value = [double]"137.0323479298667"
$ws.Range("R12").Value = [double]"1233.2911313688"
$ws.Range("S12").Value = [double]"0.0006085078130499278"
$ws.Range("T12").Value = [double]"0.0006085078130499278"
$ws.Range("G13").Value = [double]"411.37678"
$ws.Range("H13").Value = [double]"1234.13034"
$ws.Range("I13").Value = [double]"0.2862480573072345"
$ws.Range("J13").Value = [double]"0.2862480573072345"
$ws.Range("O13").Value = [double]"0.0008775937418887864"
$ws.Range("P13").Value = [double]"0.0008775937418887864"
$ws.Range("Q13").Value = [double]"56.57088927848"
$ws.Range("R13").Value = [double]"509.13800350632"
$ws.Range("S13").Value = [double]"0.0002512095037206516"
$ws.Range("T13").Value = [double]"0.0002512095037206516"
$ws.Range("G14").Value = [double]"411.37678"
$ws.Range("H14").Value = [double]"1234.13034"
$ws.Range("I14").Value = [double]"0.2862480573072345"
$ws.Range("J14").Value = [double]"0.2862480573072345"
$ws.Range("M14").Value = [double]"91.40156066666667"
$ws.Range("N14").Value = [double]"274.204682"
$ws.Range("O14").Value = [double]"0.5833025803538128"
$ws.Range("P14").Value = [double]"0.5833025803538128"
$ws.Range("Q14").Value = [double]"37600.47971402798"
$ws.Range("R14").Value = [double]"338404.3174262518"
$ws.Range("S14").Value = [double]"0.1669692304485759"
$ws.Range("T14").Value = [double]"0.1669692304485759"
$ws.Range("G15").Value = [double]"411.37678"
$ws.Range("H15").Value = [double]"1234.13034"
$ws.Range("I15").Value = [double]"0.2862480573072345"
$ws.Range("J15").Value = [double]"0.2862480573072345"
$ws.Range("M15").Value = [double]"0.5759770000000001"
$ws.Range("N15").Value = [double]"1.727931"
$ws.Range("O15").Value = [double]"0.00367574544541637"
$ws.Range("P15").Value = [double]"0.00367574544541637"
$ws.Range("Q15").Value = [double]"236.94356361406"
$ws.Range("R15").Value = [double]"2132.49207252654"
$ws.Range("S15").Value = [double]"0.001052174992906351"
$ws.Range("T15").Value = [double]"0.001052174992906351"
$ws.Range("G16").Value = [double]"411.37678"
$ws.Range("H16").Value = [double]"1234.13034"
$ws.Range("I16").Value = [double]"0.2862480573072345"
$ws.Range("J16").Value = [double]"0.2862480573072345"
$ws.Range("M16").Value = [double]"64.24849033333334"
$ws.Range("N16").Value = [double]"192.745471"
$ws.Range("O16").Value = [double]"0.4100182745450386"
$ws.Range("P16").Value = [double]"0.4100182745450385"
$ws.Range("Q16").Value = [double]"26430.33707318779"
$ws.Range("R16").Value = [double]"237873.0336586901"
$ws.Range("S16").Value = [double]"0.1173669345489816"
$ws.Range("T16").Value = [double]"0.1173669345489816"
$ws.Range("G17").Value = [double]"173.2560603333334"
$ws.Range("H17").Value = [double]"519.768181"
$ws.Range("I17").Value = [double]"0.12055666021578"
$ws.Range("J17").Value = [double]"0.12055666021578"
$ws.Range("M17").Value = [double]"0.3331066666666667"
$ws.Range("N17").Value = [double]"0.99932"
$ws.Range("O17").Value = [double]"0.002125805913843485"
$ws.Range("P17").Value = [double]"0.002125805913843485"
$ws.Range("Q17").Value = [double]"57.71274873743556"
$ws.Range("R17").Value = [double]"519.41473863692"
$ws.Range("S17").Value = [double]"0.0002562800612399247"
$ws.Range("T17").Value = [double]"0.0002562800612399247"
$ws.Range("G18").Value = [double]"173.2560603333334"
$ws.Range("H18").Value = [double]"519.768181"
$ws.Range("I18").Value = [double]"0.12055666021578"
$ws.Range("J18").Value = [double]"0.12055666021578"
$ws.Range("O18").Value = [double]"0.0008775937418887864"
$ws.Range("P18").Value = [double]"0.0008775937418887864"
$ws.Range("Q18").Value = [double]"23.82548039279867"
$ws.Range("R18").Value = [double]"214.429323535188"
$ws.Range("S18").Value = [double]"0.0001057997705483813"
$ws.Range("T18").Value = [double]"0.0001057997705483813"
$ws.Range("G19").Value = [double]"173.2560603333334"
$ws.Range("H19").Value = [double]"519.768181"
$ws.Range("I19").Value = [double]"0.12055666021578"
$ws.Range("J19").Value = [double]"0.12055666021578"
$ws.Range("M19").Value = [double]"91.40156066666667"
$ws.Range("N19").Value = [double]"274.204682"
$ws.Range("O19").Value = [double]"0.5833025803538128"
$ws.Range("P19").Value = [double]"0.5833025803538128"
$ws.Range("Q19").Value = [double]"15835.87430942483"
$ws.Range("R19").Value = [double]"142522.8687848235"
$ws.Range("S19").Value = [double]"0.0703210109827023"
$ws.Range("T19").Value = [double]"0.07032101098270231"
$ws.Range("G20").Value = [double]"173.2560603333334"
$ws.Range("H20").Value = [double]"519.768181"
$ws.Range("I20").Value = [double]"0.12055666021578"
$ws.Range("J20").Value = [double]"0.12055666021578"
$ws.Range("M20").Value = [double]"0.5759770000000001"
$ws.Range("N20").Value = [double]"1.727931"
$ws.Range("O20").Value = [double]"0.00367574544541637"
$ws.Range("P20").Value = [double]"0.00367574544541637"
$ws.Range("Q20").Value = [double]"99.79150586261235"
$ws.Range("R20").Value = [double]"898.1235527635112"
$ws.Range("S20").Value = [double]"0.0004431355947027622"
$ws.Range("T20").Value = [double]"0.0004431355947027622"
$ws.Range("G21").Value = [double]"173.2560603333334"
$ws.Range("H21").Value = [double]"519.768181"
$ws.Range("I21").Value = [double]"0.12055666021578"
$ws.Range("J21").Value = [double]"0.12055666021578"
$ws.Range("M21").Value = [double]"64.24849033333334"
$ws.Range("N21").Value = [double]"192.745471"
$ws.Range("O21").Value = [double]"0.4100182745450386"
$ws.Range("P21").Value = [double]"0.4100182745450385"
$ws.Range("Q21").Value = [double]"11131.44031751758"
$ws.Range("R21").Value = [double]"100182.9628576583"
$ws.Range("S21").Value = [double]"0.04943043380658661"
$ws.Range("T21").Value = [double]"0.04943043380658661"
$ws.Range("G22").Value = [double]"794.4973246666667"
$ws.Range("H22").Value = [double]"2383.491974"
$ws.Range("I22").Value = [double]"0.5528345953838922"
$ws.Range("J22").Value = [double]"0.5528345953838923"
$ws.Range("M22").Value = [double]"0.3331066666666667"
$ws.Range("N22").Value = [double]"0.99932"
$ws.Range("O22").Value = [double]"0.002125805913843485"
$ws.Range("P22").Value = [double]"0.002125805913843485"
$ws.Range("Q22").Value = [double]"264.6523554952978"
$ws.Range("R22").Value = [double]"2381.87119945768"
$ws.Range("S22").Value = [double]"0.001175219052244348"
$ws.Range("T22").Value = [double]"0.001175219052244349"
$ws.Range("G23").Value = [double]"794.4973246666667"
$ws.Range("H23").Value = [double]"2383.491974"
$ws.Range("I23").Value = [double]"0.5528345953838922"
$ws.Range("J23").Value = [double]"0.5528345953838923"
$ws.Range("O23").Value = [double]"0.0008775937418887864"
$ws.Range("P23").Value = [double]"0.0008775937418887864"
$ws.Range("Q23").Value = [double]"109.2560940988613"
$ws.Range("R23").Value = [double]"983.3048468897521"
$ws.Range("S23").Value = [double]"0.0004851641812085232"
$ws.Range("T23").Value = [double]"0.0004851641812085233"
$ws.Range("G24").Value = [double]"794.4973246666667"
$ws.Range("H24").Value = [double]"2383.491974"
$ws.Range("I24").Value = [double]"0.5528345953838922"
$ws.Range("J24").Value = [double]"0.5528345953838923"
$ws.Range("M24").Value = [double]"91.40156066666667"
$ws.Range("N24").Value = [double]"274.204682"
$ws.Range("O24").Value = [double]"0.5833025803538128"
$ws.Range("P24").Value = [double]"0.5833025803538128"
$ws.Range("Q24").Value = [double]"72618.2954200247"
$ws.Range("R24").Value = [double]"653564.6587802223"
$ws.Range("S24").Value = [double]"0.3224698459962804"
$ws.Range("T24").Value = [double]"0.3224698459962804"
$ws.Range("G25").Value = [double]"794.4973246666667"
$ws.Range("H25").Value = [double]"2383.491974"
$ws.Range("I25").Value = [double]"0.5528345953838922"
$ws.Range("J25").Value = [double]"0.5528345953838923"
$ws.Range("M25").Value = [double]"0.5759770000000001"
$ws.Range("N25").Value = [double]"1.727931"
$ws.Range("O25").Value = [double]"0.00367574544541637"
$ws.Range("P25").Value = [double]"0.00367574544541637"
$ws.Range("Q25").Value = [double]"457.6121855695328"
$ws.Range("R25").Value = [double]"4118.509670125794"
$ws.Range("S25").Value = [double]"0.002032079246050944"
$ws.Range("T25").Value = [double]"0.002032079246050944"
$ws.Range("G26").Value = [double]"794.4973246666667"
$ws.Range("H26").Value = [double]"2383.491974"
$ws.Range("I26").Value = [double]"0.5528345953838922"
$ws.Range("J26").Value = [double]"0.5528345953838923"
$ws.Range("M26").Value = [double]"64.24849033333334"
$ws.Range("N26").Value = [double]"192.745471"
$ws.Range("O26").Value = [double]"0.4100182745450386"
$ws.Range("P26").Value = [double]"0.4100182745450385"
$ws.Range("Q26").Value = [double]"51045.25368370553"
$ws.Range("R26").Value = [double]"459407.2831533498"
$ws.Range("S26").Value = [double]"0.2266722869081081"
$ws.Range("T26").Value = [double]"0.2266722869081081"
